# The upstream change (commit "Fixed POI packaging and upgraded to POI
# 3.15") only re-serialized this fixture's OOXML after a library upgrade:
# every hunk in the diff is a pure XML-attribute re-ordering (e.g.
# <w:pgSz w:w="11906" w:h="16838"/> -> <w:pgSz w:h="16838" w:w="11906"/>,
# attributes on <w:latentStyles>, <w:rFonts>, <w:lang>, <w:style>, ...
# all alphabetized) together with namespace-declaration reordering on the
# root <w:document> element. No text, formatting, style value, page
# geometry, language, or any other document property actually changed -
# every attribute keeps the exact same value, just written in a
# different order by the newer POI version.
#
# There is therefore no content edit to make through the Word object
# model: every value already matches what the commit produced, so the
# correct automation here is to leave the document exactly as-is rather
# than touch properties (which would only risk perturbing real content)
# merely to try to influence attribute-serialization order, something
# the object model does not expose control over in the first place.
$d = $word.ActiveDocument
